$wb = $excel.ActiveWorkbook

# Add the new "SE Results" worksheet after the last existing sheet (U238)
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "SE Results"

# Column widths
$ws.Columns.Item(1).ColumnWidth = 21.83203125
$ws.Columns.Item(2).ColumnWidth = 22.83203125
$ws.Columns.Item(3).ColumnWidth = 25.5
$ws.Columns.Item(4).ColumnWidth = 24.6640625

# Header row
$ws.Range("A1").Value = "isotopes"
$ws.Range("B1").Value = "232-Th (for 3e-14 ppm)"
$ws.Range("C1").Value = "scaled"
$ws.Range("D1").Value = "overall contamination (ppm)"

# Row 2
$ws.Range("A2").Value = "902320-Th"
$ws.Range("B2").Value = 0.00000000000003
$ws.Range("D2").Value = 0.00000000000003
$ws.Range("C2").Formula = "=B2*(`$D`$2/0.00000000000003)"

# Row 3
$ws.Range("A3").Value = "902280-Th"
$ws.Range("B3").Formula = "=4.07056E-24"
$ws.Range("C3").Formula = "=B3*(`$D`$2/0.00000000000003)"

# Row 4
$ws.Range("A4").Value = "882240-Ra"
$ws.Range("B4").Formula = "=2.11268E-26"
$ws.Range("C4").Formula = "=B4*(`$D`$2/0.00000000000003)"

# Row 5
$ws.Range("A5").Value = "862200-Rn"
$ws.Range("B5").Value = 3.7546900000000003E-30
$ws.Range("C5").Formula = "=B5*(`$D`$2/0.00000000000003)"

# Row 6
$ws.Range("A6").Value = "842160-Po"
$ws.Range("B6").Value = 9.8039200000000002E-33
$ws.Range("C6").Formula = "=B6*(`$D`$2/0.00000000000003)"

# Row 7
$ws.Range("A7").Value = "842120-Po"
$ws.Range("B7").Value = 1.2898000000000001E-38
$ws.Range("C7").Formula = "=B7*(`$D`$2/0.00000000000003)"

# Row 8
$ws.Range("A8").Value = "832120-Bi"
$ws.Range("B8").Formula = "=2.45902E-28"
$ws.Range("C8").Formula = "=B8*(`$D`$2/0.00000000000003)"

# Number formats - scientific notation for B/C columns (rows 2-8) and D2
$ws.Range("B2:C8").NumberFormat = "0.000000E+00"
$ws.Range("D2").NumberFormat = "0.000000E+00"

# Fill color for header row (B1:C1) - theme accent6, lighter 40%
$ws.Range("B1:C1").Interior.Color = RGB(169, 209, 142)

# Borders around header + data cells that have fills
$ws.Range("A1:D1").Borders.LineStyle = 1
$ws.Range("B2:C8").Borders.LineStyle = 1

# Selection & view
$ws.Range("D17").Select()

$wb.Worksheets.Item(3).Tab.Select()
$ws.Tab.Select()

$wb.Save()
